$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), P (Precio $/Kg) for rows 2-38,
# and append a brand-new row 39 with the same static fields as the others.

$ws.Cells.Item(2, 4).Value = 44249
$ws.Cells.Item(2, 10).Value = 200
$ws.Cells.Item(2, 11).Value = 900
$ws.Cells.Item(2, 12).Value = 1000
$ws.Cells.Item(2, 13).Value = 950
$ws.Cells.Item(2, 16).Value = 950

$ws.Cells.Item(3, 4).Value = 44292
$ws.Cells.Item(3, 10).Value = 250
$ws.Cells.Item(3, 11).Value = 1800
$ws.Cells.Item(3, 12).Value = 2000
$ws.Cells.Item(3, 13).Value = 1900
$ws.Cells.Item(3, 16).Value = 1900

$ws.Cells.Item(4, 4).Value = 44442
$ws.Cells.Item(4, 10).Value = 240
$ws.Cells.Item(4, 11).Value = 2300
$ws.Cells.Item(4, 12).Value = 2500
$ws.Cells.Item(4, 13).Value = 2400
$ws.Cells.Item(4, 16).Value = 2400

$ws.Cells.Item(5, 4).Value = 44349
$ws.Cells.Item(5, 10).Value = 250
$ws.Cells.Item(5, 11).Value = 2800
$ws.Cells.Item(5, 12).Value = 3000
$ws.Cells.Item(5, 13).Value = 2900
$ws.Cells.Item(5, 16).Value = 2900

$ws.Cells.Item(6, 4).Value = 44474
$ws.Cells.Item(6, 10).Value = 250
$ws.Cells.Item(6, 11).Value = 2000
$ws.Cells.Item(6, 12).Value = 2500
$ws.Cells.Item(6, 13).Value = 2250
$ws.Cells.Item(6, 16).Value = 2250

$ws.Cells.Item(7, 4).Value = 44532
$ws.Cells.Item(7, 10).Value = 300
$ws.Cells.Item(7, 11).Value = 1000
$ws.Cells.Item(7, 12).Value = 1200
$ws.Cells.Item(7, 13).Value = 1100
$ws.Cells.Item(7, 16).Value = 1100

$ws.Cells.Item(8, 4).Value = 44539
$ws.Cells.Item(8, 10).Value = 300
$ws.Cells.Item(8, 11).Value = 900
$ws.Cells.Item(8, 12).Value = 1000
$ws.Cells.Item(8, 13).Value = 950
$ws.Cells.Item(8, 16).Value = 950

$ws.Cells.Item(9, 4).Value = 44274
$ws.Cells.Item(9, 10).Value = 250
$ws.Cells.Item(9, 11).Value = 1000
$ws.Cells.Item(9, 12).Value = 1200
$ws.Cells.Item(9, 13).Value = 1100
$ws.Cells.Item(9, 16).Value = 1100

$ws.Cells.Item(10, 4).Value = 44326
$ws.Cells.Item(10, 10).Value = 200
$ws.Cells.Item(10, 11).Value = 2700
$ws.Cells.Item(10, 12).Value = 2800
$ws.Cells.Item(10, 13).Value = 2750
$ws.Cells.Item(10, 16).Value = 2750

$ws.Cells.Item(11, 4).Value = 44494
$ws.Cells.Item(11, 10).Value = 200
$ws.Cells.Item(11, 11).Value = 2400
$ws.Cells.Item(11, 12).Value = 2500
$ws.Cells.Item(11, 13).Value = 2450
$ws.Cells.Item(11, 16).Value = 2450

$ws.Cells.Item(12, 4).Value = 44628
$ws.Cells.Item(12, 10).Value = 250
$ws.Cells.Item(12, 11).Value = 2500
$ws.Cells.Item(12, 12).Value = 3000
$ws.Cells.Item(12, 13).Value = 2750
$ws.Cells.Item(12, 16).Value = 2750

$ws.Cells.Item(13, 4).Value = 44280
$ws.Cells.Item(13, 10).Value = 250
$ws.Cells.Item(13, 11).Value = 1400
$ws.Cells.Item(13, 12).Value = 1500
$ws.Cells.Item(13, 13).Value = 1450
$ws.Cells.Item(13, 16).Value = 1450

$ws.Cells.Item(14, 4).Value = 44523
$ws.Cells.Item(14, 10).Value = 250
$ws.Cells.Item(14, 11).Value = 1400
$ws.Cells.Item(14, 12).Value = 1500
$ws.Cells.Item(14, 13).Value = 1450
$ws.Cells.Item(14, 16).Value = 1450

$ws.Cells.Item(15, 4).Value = 44432
$ws.Cells.Item(15, 10).Value = 300
$ws.Cells.Item(15, 11).Value = 2300
$ws.Cells.Item(15, 12).Value = 2500
$ws.Cells.Item(15, 13).Value = 2400
$ws.Cells.Item(15, 16).Value = 2400

$ws.Cells.Item(16, 4).Value = 44540
$ws.Cells.Item(16, 10).Value = 200
$ws.Cells.Item(16, 11).Value = 900
$ws.Cells.Item(16, 12).Value = 1000
$ws.Cells.Item(16, 13).Value = 950
$ws.Cells.Item(16, 16).Value = 950

$ws.Cells.Item(17, 4).Value = 44659
$ws.Cells.Item(17, 10).Value = 250
$ws.Cells.Item(17, 11).Value = 950
$ws.Cells.Item(17, 12).Value = 1000
$ws.Cells.Item(17, 13).Value = 975
$ws.Cells.Item(17, 16).Value = 975

$ws.Cells.Item(18, 4).Value = 44376
$ws.Cells.Item(18, 10).Value = 270
$ws.Cells.Item(18, 11).Value = 2400
$ws.Cells.Item(18, 12).Value = 2500
$ws.Cells.Item(18, 13).Value = 2437
$ws.Cells.Item(18, 16).Value = 2437

$ws.Cells.Item(19, 4).Value = 44616
$ws.Cells.Item(19, 10).Value = 200
$ws.Cells.Item(19, 11).Value = 2500
$ws.Cells.Item(19, 12).Value = 3000
$ws.Cells.Item(19, 13).Value = 2750
$ws.Cells.Item(19, 16).Value = 2750

$ws.Cells.Item(20, 4).Value = 44571
$ws.Cells.Item(20, 10).Value = 250
$ws.Cells.Item(20, 11).Value = 900
$ws.Cells.Item(20, 12).Value = 1000
$ws.Cells.Item(20, 13).Value = 950
$ws.Cells.Item(20, 16).Value = 950

$ws.Cells.Item(21, 4).Value = 44536
$ws.Cells.Item(21, 10).Value = 250
$ws.Cells.Item(21, 11).Value = 900
$ws.Cells.Item(21, 12).Value = 1000
$ws.Cells.Item(21, 13).Value = 950
$ws.Cells.Item(21, 16).Value = 950

$ws.Cells.Item(22, 4).Value = 44417
$ws.Cells.Item(22, 10).Value = 250
$ws.Cells.Item(22, 11).Value = 4000
$ws.Cells.Item(22, 12).Value = 4500
$ws.Cells.Item(22, 13).Value = 4250
$ws.Cells.Item(22, 16).Value = 4250

$ws.Cells.Item(23, 4).Value = 44645
$ws.Cells.Item(23, 10).Value = 300
$ws.Cells.Item(23, 11).Value = 1200
$ws.Cells.Item(23, 12).Value = 1500
$ws.Cells.Item(23, 13).Value = 1350
$ws.Cells.Item(23, 16).Value = 1350

$ws.Cells.Item(24, 4).Value = 44362
$ws.Cells.Item(24, 10).Value = 250
$ws.Cells.Item(24, 11).Value = 2800
$ws.Cells.Item(24, 12).Value = 3000
$ws.Cells.Item(24, 13).Value = 2900
$ws.Cells.Item(24, 16).Value = 2900

$ws.Cells.Item(25, 4).Value = 44302
$ws.Cells.Item(25, 10).Value = 200
$ws.Cells.Item(25, 11).Value = 900
$ws.Cells.Item(25, 12).Value = 1000
$ws.Cells.Item(25, 13).Value = 950
$ws.Cells.Item(25, 16).Value = 950

$ws.Cells.Item(26, 4).Value = 44498
$ws.Cells.Item(26, 10).Value = 270
$ws.Cells.Item(26, 11).Value = 2000
$ws.Cells.Item(26, 12).Value = 2300
$ws.Cells.Item(26, 13).Value = 2150
$ws.Cells.Item(26, 16).Value = 2150

$ws.Cells.Item(27, 4).Value = 44635
$ws.Cells.Item(27, 10).Value = 300
$ws.Cells.Item(27, 11).Value = 1900
$ws.Cells.Item(27, 12).Value = 2000
$ws.Cells.Item(27, 13).Value = 1950
$ws.Cells.Item(27, 16).Value = 1950

$ws.Cells.Item(28, 4).Value = 44365
$ws.Cells.Item(28, 10).Value = 250
$ws.Cells.Item(28, 11).Value = 2400
$ws.Cells.Item(28, 12).Value = 2500
$ws.Cells.Item(28, 13).Value = 2450
$ws.Cells.Item(28, 16).Value = 2450

$ws.Cells.Item(29, 4).Value = 44260
$ws.Cells.Item(29, 10).Value = 250
$ws.Cells.Item(29, 11).Value = 900
$ws.Cells.Item(29, 12).Value = 1000
$ws.Cells.Item(29, 13).Value = 950
$ws.Cells.Item(29, 16).Value = 950

$ws.Cells.Item(30, 4).Value = 44330
$ws.Cells.Item(30, 10).Value = 250
$ws.Cells.Item(30, 11).Value = 2800
$ws.Cells.Item(30, 12).Value = 3000
$ws.Cells.Item(30, 13).Value = 2900
$ws.Cells.Item(30, 16).Value = 2900

$ws.Cells.Item(31, 4).Value = 44250
$ws.Cells.Item(31, 10).Value = 250
$ws.Cells.Item(31, 11).Value = 1000
$ws.Cells.Item(31, 12).Value = 1200
$ws.Cells.Item(31, 13).Value = 1100
$ws.Cells.Item(31, 16).Value = 1100

$ws.Cells.Item(32, 4).Value = 44727
$ws.Cells.Item(32, 10).Value = 270
$ws.Cells.Item(32, 11).Value = 1500
$ws.Cells.Item(32, 12).Value = 2000
$ws.Cells.Item(32, 13).Value = 1750
$ws.Cells.Item(32, 16).Value = 1750

$ws.Cells.Item(33, 4).Value = 44699
$ws.Cells.Item(33, 10).Value = 300
$ws.Cells.Item(33, 11).Value = 2000
$ws.Cells.Item(33, 12).Value = 2500
$ws.Cells.Item(33, 13).Value = 2250
$ws.Cells.Item(33, 16).Value = 2250

$ws.Cells.Item(34, 4).Value = 44664
$ws.Cells.Item(34, 10).Value = 250
$ws.Cells.Item(34, 11).Value = 1300
$ws.Cells.Item(34, 12).Value = 1500
$ws.Cells.Item(34, 13).Value = 1400
$ws.Cells.Item(34, 16).Value = 1400

$ws.Cells.Item(35, 4).Value = 44435
$ws.Cells.Item(35, 10).Value = 300
$ws.Cells.Item(35, 11).Value = 2300
$ws.Cells.Item(35, 12).Value = 2500
$ws.Cells.Item(35, 13).Value = 2400
$ws.Cells.Item(35, 16).Value = 2400

$ws.Cells.Item(36, 4).Value = 44669
$ws.Cells.Item(36, 10).Value = 300
$ws.Cells.Item(36, 11).Value = 950
$ws.Cells.Item(36, 12).Value = 1000
$ws.Cells.Item(36, 13).Value = 975
$ws.Cells.Item(36, 16).Value = 975

$ws.Cells.Item(37, 4).Value = 44603
$ws.Cells.Item(37, 10).Value = 250
$ws.Cells.Item(37, 11).Value = 2500
$ws.Cells.Item(37, 12).Value = 3000
$ws.Cells.Item(37, 13).Value = 2750
$ws.Cells.Item(37, 16).Value = 2750

$ws.Cells.Item(38, 4).Value = 44224
$ws.Cells.Item(38, 10).Value = 200
$ws.Cells.Item(38, 11).Value = 750
$ws.Cells.Item(38, 12).Value = 800
$ws.Cells.Item(38, 13).Value = 775
$ws.Cells.Item(38, 16).Value = 775

$ws.Cells.Item(39, 4).Value = 44313
$ws.Cells.Item(39, 10).Value = 250
$ws.Cells.Item(39, 11).Value = 900
$ws.Cells.Item(39, 12).Value = 1000
$ws.Cells.Item(39, 13).Value = 950
$ws.Cells.Item(39, 16).Value = 950

# Row 39 is brand new: give its Fecha cell the same date number format as
# the other rows in column D (style s="2", numFmtId 165) before/while
# setting its value.
$ws.Cells.Item(39, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 39 needs the remaining (previously unused) static columns too
$ws.Cells.Item(39, 1).Value = 1
$ws.Cells.Item(39, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(39, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(39, 5).Value = 15
$ws.Cells.Item(39, 6).Value = 100112052
$ws.Cells.Item(39, 7).Value = "Albahaca"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 14).Value = "$/paquete"
$ws.Cells.Item(39, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(39, 17).Value = 1
$ws.Cells.Item(39, 18).Value = "Hortaliza"
